# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "45.015.34"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.28%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.264.05"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.32%  "

# Row 4
$ws.Range("E4").Value = "  -0.89%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.79"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.69%  "

# Row 7
$ws.Range("E7").Value = "  -0.37%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.62%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.76%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.06"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.42%  "

# Row 11
$ws.Range("E11").Value = "  -0.28%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.20"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.42%  "

# Row 13
$ws.Range("E13").Value = "  -1.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.609.98"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.64%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.268.41"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.96%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.806"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.37%  "

# Row 17
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.54"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.95%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.932.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.71%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.53"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +15.72%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0915"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.19%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.99"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.53%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.32"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.92%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.81"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.07%  "

# Row 24
$ws.Range("E24").Value = "  -0.75%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.66%  "

# Row 26
$ws.Range("E26").Value = "  -2.03%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.57"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +12.37%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.23"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.24%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.53"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.57"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.93%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.35"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.16%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.51"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -4.91%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0791"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.75%  "

# Row 34
$ws.Range("E34").Value = "  -2.13%  "

# Row 35
$ws.Range("E35").Value = "  -0.67%  "

# Row 36
$ws.Range("E36").Value = "  -6.61%  "

# Row 37
$ws.Range("E37").Value = "  -4.06%  "

# Row 38
$ws.Range("E38").Value = "  -4.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0311"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.02%  "

# Row 41
$ws.Range("E41").Value = "  -3.53%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.56"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.16%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.90%  "

# Row 44
$ws.Range("E44").Value = "  +12.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.772.62"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.49%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.191"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.14%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "69.91"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.84%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "76.01"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.16%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.06"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.41%  "

# Row 50
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.70"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.86"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.56%  "
